$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 7 (pushes the old rows 7..19 down to 8..20,
# and extends the SUM/COUNTA ranges in the TOTAL row automatically).
$ws.Rows.Item(7).Insert()

# The freshly inserted row doesn't reliably pick up the same border/fill
# formatting as its neighbours, so copy the (now shifted-down) row 8
# formatting - which used to be row 7 before the insert - back onto the
# new row 7.
$ws.Range("C8:I8").Copy()
$ws.Range("C7:I7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new tile entry: name in C, appearance counts D:G left
# blank (untested), "Juntas Testadas?" (H) left blank, and the "x"
# marker in I that every data row carries.
$ws.Range("C7").Value = "BP_Curva_Baixo_01"
$ws.Range("I7").Value = "x"

# Re-apply the TOTAL row's D:F formula as one block so it is re-grouped
# into a shared formula the same way it was before the insert (now
# covering the extended D20:F20 / row 4-19 range).
$ws.Range("D20:F20").Formula = "=SUM(D4:D19)"

# Reflect where the author was last working.
$ws.Range("C7").Select()
